$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column B (cells B1, B2 removed entirely)
$ws.Range("B1:B2").Clear()

# A1 keeps its existing "yyyy-mm-dd" style, just change the value
$ws.Range("A1").Value = 45306

# A2 keeps the same value but switches to the "m/d/yy" number format (style index 2)
$ws.Range("A2").Value = 45306
$ws.Range("A2").NumberFormat = "m/d/yy;@"

# New row A3 uses the "yyyy-mm-dd" number format (style index 1), same as A1
$ws.Range("A3").Value = 45306
$ws.Range("A3").NumberFormat = "yyyy\-mm\-dd;@"

# Clear the stale selection left over on B3
$ws.Range("A1").Select()
